$d = $word.ActiveDocument
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Import dump*") {
        # Insert 3 paragraph breaks right after the current (still plain-formatted) paragraph
        $r = $p.Range
        $r.Collapse(1)
        $r.InsertParagraphAfter()
        $p2 = $p.Next()
        $r2 = $p2.Range
        $r2.Collapse(1)
        $r2.InsertParagraphAfter()
        $p3 = $p2.Next()
        $r3 = $p3.Range
        $r3.Collapse(1)
        $r3.InsertParagraphAfter()
        $p4 = $p3.Next()

        # Now set text+format on paragraph 1 (Import latest dump...)
        $p.Format.SpaceAfter = 0
        $p.Range.Text = "Import latest dump located at "
        $p.Range.Font.Name = "Verdana"
        $p.Range.Font.Size = 10

        # paragraph 2 (Oracle)
        $p2.Format.SpaceAfter = 0
        $p2.Range.Text = "Oracle: https://ncisvn.nci.nih.gov/svn/catissue_persistent/caTissue Database Dump/v2.0/Oracle"
        $p2.Range.Font.Name = "Verdana"
        $p2.Range.Font.Size = 10

        # paragraph 3 (MySQL)
        $p3.Format.SpaceAfter = 0
        $p3.Range.Text = "MySQL: https://ncisvn.nci.nih.gov/svn/catissue_persistent/caTissue Database Dump/v2.0/MySQL and deploy application."
        $p3.Range.Font.Name = "Verdana"
        $p3.Range.Font.Size = 10

        # paragraph 4 (blank, bold) -- left untouched except bold
        $p4.Range.Font.Bold = 1

        Write-Host "done"
    }
}
Write-Host $d.Paragraphs.Count
